$wb = $excel.ActiveWorkbook

# --- Rename "Sheet2" -> "webshop" ------------------------------------------
$wsShop = $wb.Worksheets.Item("Sheet2")
$wsShop.Name = "webshop"

# --- Populate the new "webshop" sheet with header + test-data row ----------
$wsShop.Range("A1").Value = "UserName"
$wsShop.Range("B1").Value = "Password"

# Write the password first, then the email - matches the shared-string order
# the workbook ends up with (password allocated before the email address).
$wsShop.Range("B2").Value = "Sample@1234567"
$wsShop.Range("A2").Value = "vindhya1.tech@gmail.com"

# Column widths (roughly matching the authored widths; engine quantizes to
# 1/6-character steps so this is the closest achievable approximation).
$wsShop.Columns.Item(1).ColumnWidth = 21.0
$wsShop.Columns.Item(2).ColumnWidth = 16.666666666666668

# Hyperlinks on the credential row - Excel auto-links a typed e-mail address;
# both cells end up pointing at the same mailto: address.
$wsShop.Hyperlinks.Add($wsShop.Range("A2"), "mailto:vindhya1.tech@gmail.com")
$wsShop.Hyperlinks.Add($wsShop.Range("B2"), "mailto:vindhya1.tech@gmail.com")

# --- Make "webshop" the active/selected sheet -------------------------------
$wsShop.Activate()
$wsShop.Range("A2").Select()

Write-Output "done"
